$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# The merged header cells B1:D1 (and E1:G1 on the other sheet) are drawn with
# a thin box border. B1/E1 already carry that border via the existing bold
# style; give the filler cells C1/D1 (and F1/G1) the matching interior border
# pieces: a top+bottom edge for the middle cell of the box, and a
# top+right+bottom edge for the cell that closes the box on the right.
#
# Build it as "full box, then drop the sides that shouldn't be there" so the
# engine lands on the pre-existing border definitions (4 = top+bottom,
# 5 = top+right+bottom) instead of minting new ones.
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.LineStyle = 1
$c1.Borders.Item(7).LineStyle = -4142
$c1.Borders.Item(10).LineStyle = -4142

$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.LineStyle = 1
$d1.Borders.Item(7).LineStyle = -4142

# Every other filler cell needs one of those same two formats, so copy the
# already-built formatting across rather than re-deriving it (keeps the
# style table from growing beyond the two new entries that are needed).
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("F1").PasteSpecial(-4122)

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Anonymize: rename the "fedcore" column headers to "approach".
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 was a stray empty inline-string cell; drop it entirely.
$ws2.Range("G5").ClearContents()
